# Scheduled runner update: refresh market-price derived columns (H,I,J,K,L,M,N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR Leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12 (Leve Item ID 5515)
$ws.Range("H12").Value = 383.33334
$ws.Range("I12").Value = 175
$ws.Range("K12").Value = 175
$ws.Range("M12").Value = -5
# Row 31 (Leve Item ID 4576)
$ws.Range("H31").Value = 246.14285
$ws.Range("I31").Value = 246.14285
$ws.Range("K31").Value = 738.4285500000001
$ws.Range("M31").Value = -508.4285500000001
# Row 39 (Leve Item ID 4603)
$ws.Range("H39").Value = 253.16667
$ws.Range("I39").Value = 153.5
$ws.Range("K39").Value = 460.5
$ws.Range("M39").Value = -164.5
# Row 41 (Leve Item ID 5478)
$ws.Range("H41").Value = 398.85715
$ws.Range("I41").Value = 398.85715
$ws.Range("K41").Value = 398.85715
$ws.Range("M41").Value = 41.14285000000001
# Row 82 (Leve Item ID 12623)
$ws.Range("H82").Value = 683
$ws.Range("I82").Value = 683
$ws.Range("K82").Value = 2049
$ws.Range("M82").Value = -1643
# Row 85 (Leve Item ID 12623)
$ws.Range("H85").Value = 683
$ws.Range("I85").Value = 683
$ws.Range("K85").Value = 2049
$ws.Range("M85").Value = -645
# Row 116 (Leve Item ID 27778)
$ws.Range("H116").Value = 4897.3335
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 4897.3335
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 4897.3335
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -11781.3335
$ws = $wb.Worksheets.Item("ARM")
# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 2022.5
$ws.Range("J45").Value = 2000
$ws.Range("L45").Value = 2000
$ws.Range("N45").Value = -2754
# Row 53 (Leve Item ID 3623)
$ws.Range("H53").Value = 12999
$ws.Range("I53").Value = 12999
$ws.Range("K53").Value = 12999
$ws.Range("M53").Value = -12317
# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 4757
$ws.Range("J61").Value = 5014
$ws.Range("L61").Value = 5014
$ws.Range("N61").Value = -5438
# Row 63 (Leve Item ID 12528)
$ws.Range("H63").Value = 5008.1665
$ws.Range("I63").Value = 5008.1665
$ws.Range("K63").Value = 5008.1665
$ws.Range("M63").Value = -4322.1665
# Row 66 (Leve Item ID 12528)
$ws.Range("H66").Value = 5008.1665
$ws.Range("I66").Value = 5008.1665
$ws.Range("K66").Value = 25040.8325
$ws.Range("M66").Value = -21608.8325
# Row 95 (Leve Item ID 18204)
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
# Row 98 (Leve Item ID 18371)
$ws.Range("H98").Value = 57098.6
$ws.Range("J98").Value = 57098.6
$ws.Range("L98").Value = 57098.6
$ws.Range("N98").Value = -63088.6
# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 4757
$ws.Range("J136").Value = 5014
$ws.Range("L136").Value = 15042
$ws.Range("N136").Value = -20142
$ws = $wb.Worksheets.Item("BSM")
# Row 22 (Leve Item ID 5092)
$ws.Range("H22").Value = 9998
$ws.Range("I22").Value = 9998
$ws.Range("K22").Value = 9998
$ws.Range("M22").Value = -9825
# Row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 881.6
$ws.Range("I86").Value = 881.6
$ws.Range("K86").Value = 881.6
$ws.Range("M86").Value = 241.4
# Row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 881.6
$ws.Range("I89").Value = 881.6
$ws.Range("K89").Value = 4408
$ws.Range("M89").Value = 1208
$ws = $wb.Worksheets.Item("CRP")
# Row 26 (Leve Item ID 2004)
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
# Row 55 (Leve Item ID 1855)
$ws.Range("H55").Value = 10579.167
$ws.Range("I55").Value = 4500
$ws.Range("J55").Value = 11795
$ws.Range("K55").Value = 4500
$ws.Range("L55").Value = 11795
$ws.Range("M55").Value = -4185
$ws.Range("N55").Value = -12425
# Row 58 (Leve Item ID 44021)
$ws.Range("H58").Value = 8662
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
# Row 108 (Leve Item ID 27087)
$ws.Range("H108").Value = 45000
$ws.Range("J108").Value = 40000
$ws.Range("L108").Value = 40000
$ws.Range("N108").Value = -47680
# Row 122 (Leve Item ID 36196)
$ws.Range("H122").Value = 6283.5713
$ws.Range("I122").Value = 6830.8335
$ws.Range("K122").Value = 20492.5005
$ws.Range("M122").Value = -18042.5005
# Row 136 (Leve Item ID 44021)
$ws.Range("H136").Value = 8662
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
# Row 141 (Leve Item ID 43345)
$ws.Range("H141").Value = 279368.75
$ws.Range("J141").Value = 408738.5
$ws.Range("L141").Value = 408738.5
$ws.Range("N141").Value = -419098.5
$ws = $wb.Worksheets.Item("CUL")
# Row 5 (Leve Item ID 43974)
$ws.Range("H5").Value = 1849.6666
$ws.Range("J5").Value = 1399.3334
$ws.Range("L5").Value = 4198.0002
$ws.Range("N5").Value = -4422.0002
# Row 9 (Leve Item ID 4681)
$ws.Range("H9").Value = 271.7143
$ws.Range("I9").Value = 50.666668
$ws.Range("K9").Value = 152.000004
$ws.Range("M9").Value = 71.99999600000001
# Row 34 (Leve Item ID 4749)
$ws.Range("H34").Value = 400
$ws.Range("I34").Value = 400
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1200
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1116
$ws.Range("N34").ClearContents()
# Row 39 (Leve Item ID 4712)
$ws.Range("H39").Value = 1104.6666
$ws.Range("J39").Value = 1249.5
$ws.Range("L39").Value = 3748.5
$ws.Range("N39").Value = -4336.5
# Row 51 (Leve Item ID 4646)
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("N51").ClearContents()
# Row 135 (Leve Item ID 43974)
$ws.Range("H135").Value = 1849.6666
$ws.Range("J135").Value = 1399.3334
$ws.Range("L135").Value = 12594.0006
$ws.Range("N135").Value = -17664.0006
# Row 140 (Leve Item ID 44097)
$ws.Range("H140").Value = 1376.3334
$ws.Range("I140").Value = 1376.3334
$ws.Range("K140").Value = 4129.0002
$ws.Range("M140").Value = 1050.9998
$ws = $wb.Worksheets.Item("GSM")
# Row 36 (Leve Item ID 4222)
$ws.Range("H36").Value = 14703.333
$ws.Range("J36").Value = 18055
$ws.Range("L36").Value = 18055
$ws.Range("N36").Value = -19025
# Row 45 (Leve Item ID 27225)
$ws.Range("H45").Value = 100000
$ws.Range("J45").Value = 100000
$ws.Range("L45").Value = 100000
$ws.Range("N45").Value = -101118
# Row 57 (Leve Item ID 2876)
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
# Row 63 (Leve Item ID 11048)
$ws.Range("H63").Value = 20000
$ws.Range("I63").Value = 20000
$ws.Range("K63").Value = 20000
$ws.Range("M63").Value = -19314
# Row 66 (Leve Item ID 11048)
$ws.Range("H66").Value = 20000
$ws.Range("I66").Value = 20000
$ws.Range("K66").Value = 60000
$ws.Range("M66").Value = -56568
# Row 102 (Leve Item ID 36169)
$ws.Range("H102").Value = 2204.8
$ws.Range("I102").Value = 2204.8
$ws.Range("K102").Value = 2204.8
$ws.Range("M102").Value = -582.8000000000002
# Row 104 (Leve Item ID 18666)
$ws.Range("H104").Value = 40671
$ws.Range("J104").Value = 40671
$ws.Range("L104").Value = 40671
$ws.Range("N104").Value = -47659
# Row 126 (Leve Item ID 36184)
$ws.Range("H126").Value = 2256
$ws.Range("I126").Value = 2256
$ws.Range("K126").Value = 6768
$ws.Range("M126").Value = -4298
# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
# Row 42 (Leve Item ID 4333)
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
# Row 49 (Leve Item ID 4333)
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
# Row 101 (Leve Item ID 18549)
$ws.Range("H101").Value = 14452
$ws.Range("J101").Value = 14452
$ws.Range("L101").Value = 14452
$ws.Range("N101").Value = -20942
# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
# Row 103 (Leve Item ID 18548)
$ws.Range("H103").Value = 20602
$ws.Range("J103").Value = 20602
$ws.Range("L103").Value = 20602
$ws.Range("N103").Value = -22946
# Row 105 (Leve Item ID 18710)
$ws.Range("H105").Value = 23997.5
$ws.Range("J105").Value = 23997.5
$ws.Range("L105").Value = 23997.5
$ws.Range("N105").Value = -30985.5
# Row 107 (Leve Item ID 27746)
$ws.Range("H107").Value = 1817.8334
$ws.Range("I107").Value = 1725.25
$ws.Range("J107").Value = 2003
$ws.Range("K107").Value = 5175.75
$ws.Range("L107").Value = 6009
$ws.Range("M107").Value = -3255.75
$ws.Range("N107").Value = -9849
# Row 126 (Leve Item ID 36210)
$ws.Range("H126").Value = 3575.375
$ws.Range("I126").Value = 4320.8
$ws.Range("J126").Value = 2333
$ws.Range("K126").Value = 12962.4
$ws.Range("L126").Value = 6999
$ws.Range("M126").Value = -10492.4
$ws.Range("N126").Value = -11939
# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 18000
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()
